# Insert a new data row above current row 54 (shifts rows 54:91 down to 55:92)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly price record
$ws.Cells.Item(54, 1).Value = 7
$ws.Cells.Item(54, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value = "Ñuble"
$ws.Cells.Item(54, 4).Value = 44488
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = 100112045
$ws.Cells.Item(54, 7).Value = "Zapallo"
$ws.Cells.Item(54, 8).Value = "Camote"
$ws.Cells.Item(54, 9).Value = "1a (guarda)"
$ws.Cells.Item(54, 10).Value = 120
$ws.Cells.Item(54, 11).Value = 800
$ws.Cells.Item(54, 12).Value = 900
$ws.Cells.Item(54, 13).Value = 850
$ws.Cells.Item(54, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value = 850
$ws.Cells.Item(54, 17).Value = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
